# Applies the commit "ppt and lecture notes":
#  - Slide 6 ("Purpose of Data Visualization"): bullet text
#    "Distribution in a wide range" -> "Heatmap"
#  - Slide 19 ("Homework"):
#      * give the slide its own white solid background fill
#      * split the run "Add description for what you are visualizing..."
#        so that the word "what" gets its own run in the
#        "Helvetica (light)" typeface

$p = $ppt.ActivePresentation

function Get-BodyPlaceholder($slide) {
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $candidate = $slide.Shapes.Item($j)
        if ($candidate.HasTextFrame -and $candidate.PlaceholderFormat.Type -eq 2) {
            return $candidate
        }
    }
    # fall back to the second shape, which is the body placeholder on
    # both slides touched by this edit
    return $slide.Shapes.Item(2)
}

# ---------------------------------------------------------------------
# Slide 6: "Distribution in a wide range" -> "Heatmap"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$body6 = Get-BodyPlaceholder $s6
$tr6 = $body6.TextFrame.TextRange
$paraCount6 = $tr6.Paragraphs().Count
for ($i = 1; $i -le $paraCount6; $i++) {
    $para = $tr6.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd("`r")
    if ($paraText -eq "Distribution in a wide range") {
        $para.Text = "Heatmap"
    }
}

# ---------------------------------------------------------------------
# Slide 19: white background + split "what" into its own run/font
# ---------------------------------------------------------------------
$s19 = $p.Slides.Item(19)

# Give the slide an explicit solid white background fill
$s19.Background.Fill.Solid()
$s19.Background.Fill.ForeColor.RGB = 16777215

# Locate the body placeholder and the paragraph containing the target text
$body19 = Get-BodyPlaceholder $s19
$tr19 = $body19.TextFrame.TextRange
$paraCount19 = $tr19.Paragraphs().Count
for ($i = 1; $i -le $paraCount19; $i++) {
    $para = $tr19.Paragraphs($i, 1)
    $idx = $para.Text.IndexOf("what you are visualizing")
    if ($idx -ge 0) {
        $word = $para.Characters($idx + 1, 4)
        $word.Font.Name = "Helvetica (light)"
    }
}

Write-Host "edit complete"
